$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update scenario inputs
$ws.Range("E2").Value = 71
$ws.Range("K2").Value = 550
$ws.Range("E4").Value = 0.1
$ws.Range("K4").Value = 1
$ws.Range("E12").Value = 103

# Remove the now-unused helper row (D18: =152/5)
$ws.Rows.Item(18).Delete()

# Update the selection to match the new state
$ws.Range("F19").Select()
